$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Registering PIs/PoV/PoC from Projects/Hackathon/other events ..." --
#    split the trailing run so "PoC" and "Hackathon" get their own
#    spell-check proofErr wrapping (matches Word's live spell-check pass).
# ---------------------------------------------------------------------------
$p9 = $d.Paragraphs.Item(9)
$rng9 = $p9.Range
if ($rng9.Text -like "*Registering PIs*") {
    $xml9 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="45F181AA" w14:textId="48A785F1" w:rsidR="00235519" w:rsidRDefault="00235519" w:rsidP="00235519"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Registering PIs/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PoV</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PoC</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> from Projects/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Hackathon</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/other events and is accessible for reference</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng9.InsertXML($xml9)
}

# ---------------------------------------------------------------------------
# 2) "Ideathon/Hackathon can be conducted with the help of the Portal ..." --
#    split the trailing run so "Hackathon" gets its own proofErr wrapping.
# ---------------------------------------------------------------------------
$p11 = $d.Paragraphs.Item(11)
$rng11 = $p11.Range
if ($rng11.Text -like "*Ideathon*") {
    $xml11 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5D08364A" w14:textId="06AFC01C" w:rsidR="00235519" w:rsidRDefault="00235519" w:rsidP="00235519"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="3"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Ideathon</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Hackathon</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> can be conducted with the help of the Portal and the inbuilt Framework buil</w:t></w:r><w:r w:rsidR="000A1A9D"><w:t>t</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng11.InsertXML($xml11)
}

# ---------------------------------------------------------------------------
# 3) Roles list: drop the "Operation Engineer" bullet (and the stray empty
#    paragraph that followed it), and stamp a "_GoBack" bookmark at the end
#    of the "Testers." bullet (last-edit marker Word leaves behind).
# ---------------------------------------------------------------------------
$testersPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Testers.*") {
        $testersPara = $d.Paragraphs.Item($i)
        break
    }
}

if ($testersPara -ne $null) {
    $tRng = $testersPara.Range
    $xmlT = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="378C9B84" w14:textId="7499B210" w:rsidR="00235519" w:rsidRDefault="00235519" w:rsidP="00235519"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Testers.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $tRng.InsertXML($xmlT)
}

# Now find & remove the "Operation Engineer" bullet plus the empty
# paragraph (ind left=720) that immediately follows it.
$opPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Operation Engineer*") {
        $opPara = $i
        break
    }
}

if ($opPara -ne $null) {
    $startRng = $d.Paragraphs.Item($opPara).Range
    $endRng = $d.Paragraphs.Item($opPara + 1).Range
    $delRange = $d.Range($startRng.Start, $endRng.End)
    $delRange.Delete()
}
